# Update the "人气"/view-count style numbers (column F) on several rows
# across the "展览", "演出" and "全部类型" sheets, per the source diff.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value  = 565
$wsExhibition.Range("F6").Value  = 787
$wsExhibition.Range("F7").Value  = 358
$wsExhibition.Range("F10").Value = 1116
$wsExhibition.Range("F11").Value = 568
$wsExhibition.Range("F19").Value = 513
$wsExhibition.Range("F21").Value = 459

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F5").Value = 95

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value  = 565
$wsAll.Range("F10").Value = 787
$wsAll.Range("F11").Value = 358
$wsAll.Range("F14").Value = 1116
$wsAll.Range("F15").Value = 568
$wsAll.Range("F16").Value = 95
$wsAll.Range("F31").Value = 513
$wsAll.Range("F33").Value = 459
